$wb = $excel.ActiveWorkbook

# zh-cn sheet: the Handback run that previously produced
# "2016-03-11 05:45:02" (rows 5,7,8,9,11,12,14,15,16) and
# "2016-03-11 05:45:38" (rows 10,13) now all completed at 05:46:15.
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhHandback = "2016-03-11 05:46:15"
$wsZh.Range("D5").Value = $zhHandback
$wsZh.Range("D7").Value = $zhHandback
$wsZh.Range("D8").Value = $zhHandback
$wsZh.Range("D9").Value = $zhHandback
$wsZh.Range("D10").Value = $zhHandback
$wsZh.Range("D11").Value = $zhHandback
$wsZh.Range("D12").Value = $zhHandback
$wsZh.Range("D13").Value = $zhHandback
$wsZh.Range("D14").Value = $zhHandback
$wsZh.Range("D15").Value = $zhHandback
$wsZh.Range("D16").Value = $zhHandback

# de-de sheet: same handoff/handback batch, localized for de-de,
# previously "2016-03-11 05:45:11" / "2016-03-11 05:45:46", now 05:46:24.
$wsDe = $wb.Worksheets.Item("de-de")
$deHandback = "2016-03-11 05:46:24"
$wsDe.Range("D5").Value = $deHandback
$wsDe.Range("D7").Value = $deHandback
$wsDe.Range("D8").Value = $deHandback
$wsDe.Range("D9").Value = $deHandback
$wsDe.Range("D10").Value = $deHandback
$wsDe.Range("D11").Value = $deHandback
$wsDe.Range("D12").Value = $deHandback
$wsDe.Range("D13").Value = $deHandback
$wsDe.Range("D14").Value = $deHandback
$wsDe.Range("D15").Value = $deHandback
$wsDe.Range("D16").Value = $deHandback
